$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "328.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.31%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "43.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.31%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.571"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.22%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08110"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.78%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.671"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.78%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.908"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-3.96%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "4.294"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-4.85%"

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.82%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9491"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.38%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1180"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.81%"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.49%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09640"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.78%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04109"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.31%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.1067"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.25%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001280"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.83%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005912"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-3.10%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.12%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3484"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.73%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.647"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-5.15%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.10%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2588"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.57%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04318"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.26%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001239"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.29%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004395"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.28%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "3.37%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003993"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.15%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02671"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.04%"

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.05%"

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "26.20%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007652"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.58%"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.07%"

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "1.32%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009757"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.85%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007018"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.15%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.14%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003453"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.23%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002272"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.48%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.14%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.14%"
